# Book1.xlsx -- "Added Asserter/LogicHooks to Excel compiler"
#
# The workbook originally held three near-duplicate sheets
# ("1-Schritt-Dialoge", "2-Schritt-Dialoge", "3-Schritt-Dialoge").
# The edit collapses this down to a single sheet named "Convos" (built
# from the old "3-Schritt-Dialoge" sheet), relabels the conversation
# header row as User/Bot, and augments two of the bot lines with
# scripting directives (PAUSE / BUTTONS) - wrapping those two cells so
# the extra text is visible.

$wb = $excel.ActiveWorkbook

# --- Drop the two sheets that are no longer needed ------------------------
[void]$wb.Worksheets.Item("1-Schritt-Dialoge").Delete()
[void]$wb.Worksheets.Item("2-Schritt-Dialoge").Delete()

# --- Rename the remaining sheet and make it the active one -----------------
$ws = $wb.Worksheets.Item("3-Schritt-Dialoge")
$ws.Name = "Convos"
$ws.Activate()

# --- Header row: User / Bot -------------------------------------------------
$ws.Range("A1").Value = "User"
$ws.Range("B1").Value = "Bot"

# --- Conversation body -------------------------------------------------------
$ws.Range("A2").Value  = "Ich will kündigen"
$ws.Range("B3").Value  = "Was genau?`nBUTTONS Button1|Button2|Button3"
$ws.Range("A4").Value  = "Telefonie"
$ws.Range("B5").Value  = "Möchten Sie….?`nPAUSE 1000"
$ws.Range("A6").Value  = "Nein"
$ws.Range("B7").Value  = "Hier können Sie kündigen"
$ws.Range("A9").Value  = "Kündigung"
$ws.Range("B10").Value = "Was genau?"
$ws.Range("A11").Value = "Zusatzpaket"
$ws.Range("B12").Value = "Möchten Sie….?"
$ws.Range("A13").Value = "Ja"
$ws.Range("B14").Value = "Rufen Sie an…."

# --- Wrap + taller rows for the two cells that now carry extra script lines --
$ws.Range("B3").WrapText = $true
$ws.Range("B5").WrapText = $true
$ws.Rows.Item(3).RowHeight = 28.8
$ws.Rows.Item(5).RowHeight = 28.8

# --- Selection matches the saved view in the target file -------------------
[void]$ws.Range("B6").Select()
